$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new product ("1 2 3 (ONE TWO THREE) 20 F.C.TABS.") was added as the very
# first row of the data table, pushing the existing six product rows (and the
# totals / footer rows below them) down by one row.
$ws.Rows("4:4").Insert()

# Give the freshly inserted row 4 the same cell formatting (styles, merges,
# fonts, borders, number formats) as the data rows around it by copying the
# format of the row that now sits directly below it (the old row 4, now row 5).
$ws.Range("A5:N5").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)

# Re-create the same merged-cell layout used by every other data row.
$ws.Range("B4:G4").Merge()
$ws.Range("H4:K4").Merge()
$ws.Range("L4:M4").Merge()

# Match the explicit row height used by this row position.
$ws.Rows(4).RowHeight = 24.75

# Fill in the new product's data.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "1 2 3 (ONE TWO THREE) 20 F.C.TABS."
$ws.Range("H4").Value = "2:1"
$ws.Range("L4").Value = 40
$ws.Range("N4").Value = 1

# Renumber the "م" (sequence number) column for the rows that shifted down.
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 7

# Update the grand-total cell (now on row 11) to include the new product's price.
$ws.Range("K11").Value = 515.20000000000005

# The totals row picks up a slightly different explicit height in the new layout.
$ws.Rows(11).RowHeight = 25.5
